$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Header row: "Action " -> "Action" (drop trailing space) ---
$ws.Range("A1").Value = "Action"

# --- Insert two new rows after row 3 (LoginUFT, Delay 60000) ---
$ws.Rows.Item(4).Insert()
$ws.Rows.Item(4).Insert()

$ws.Range("A4").Value = "LoginUFT"

$ws.Range("A5").Value = "Delay "
$ws.Range("B5").Value = 60000

# --- Old row4 (ModelSelect/AUTO) is now row6: leave as-is ---

# --- Old row5 (Delay/10000) is now row7: bump the delay amount to 20000 ---
$ws.Range("B7").Value = 20000

# --- Old rows 6,7,8 (SelectButton/OK, SelectButton/OK, LaunchApplication/Fault finding)
#     are now rows 8,9,10: leave as-is ---

# --- Old row9 (SelectECU/Family/SubFamily) is now row11: update the ECU detail ---
$ws.Range("B11").Value = "Engine management ECU"
$ws.Range("C11").Value = "CMM_MD1CS003"

# --- Append two new rows at the end of the script ---
$ws.Range("A12").Value = "SelectButton"
$ws.Range("B12").Value = "OK"

$ws.Range("A13").Value = "SelectMenu"
$ws.Range("B13").Value = "IDENTIFICATION"

# --- Column widths (best-effort; engine quantizes to 1/6 character units) ---
$ws.Columns.Item(1).ColumnWidth = 17.333333333333332
$ws.Columns.Item(2).ColumnWidth = 24.5
$ws.Columns.Item(3).ColumnWidth = 20.333333333333332
$ws.Columns.Item(4).ColumnWidth = 14.666666666666666
$ws.Columns.Item(5).ColumnWidth = 13.666666666666666
$ws.Columns.Item(6).ColumnWidth = 17.833333333333332

# --- Selection moves to B15 ---
$null = $ws.Range("B15").Select()
